$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records for "Choclo" / "Dulce o Americano" (Región de
# Arica y Parinacota, $/malla 70 unidades) need to be inserted into the
# chronological data table. Inserting whole rows pushes all the existing
# rows below them down, which reproduces the row-shift seen in the diff.

# --- New row inserted before current row 17 -------------------------------
$ws.Rows(17).Insert()
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "Vega Monumental Concepción"
$ws.Range("C17").Value = "Bíobío"
$ws.Range("D17").Value = 44427
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 100112024
$ws.Range("G17").Value = "Choclo"
$ws.Range("H17").Value = "Dulce o Americano"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 34000
$ws.Range("L17").Value = 35000
$ws.Range("M17").Value = 34500
$ws.Range("N17").Value = "`$/malla 70 unidades"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 493
$ws.Range("Q17").Value = 70
$ws.Range("R17").Value = "Hortaliza"

# --- New row inserted before current row 33 (after the first insert) ------
$ws.Rows(33).Insert()
$ws.Range("A33").Value = 11
$ws.Range("B33").Value = "Vega Monumental Concepción"
$ws.Range("C33").Value = "Bíobío"
$ws.Range("D33").Value = 44420
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = 100112024
$ws.Range("G33").Value = "Choclo"
$ws.Range("H33").Value = "Dulce o Americano"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 34000
$ws.Range("L33").Value = 35000
$ws.Range("M33").Value = 34500
$ws.Range("N33").Value = "`$/malla 70 unidades"
$ws.Range("O33").Value = "Región de Arica y Parinacota"
$ws.Range("P33").Value = 493
$ws.Range("Q33").Value = 70
$ws.Range("R33").Value = "Hortaliza"
